# Update "want to go" (想去人数) counts after a fresh scrape run.
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 120   # 合肥·灵能百分百ONLY2.0      119 -> 120
$ws1.Range("F8").Value = 124   # 合肥·首届进击的巨人ONLY漫展  123 -> 124
$ws1.Range("F10").Value = 6899 # 合肥·第七届环形宇宙动漫游戏嘉年华 6891 -> 6899
$ws1.Range("F12").Value = 384  # 合肥·比翼连枝国乙&代号鸢only  383 -> 384
$ws1.Range("F13").Value = 3257 # 合肥·第八届环形宇宙动漫游戏嘉年华Plus 3244 -> 3257
$ws1.Range("F15").Value = 397  # ...水千丞签售预约票           395 -> 397
$ws1.Range("F17").Value = 566  # 合肥·SSS第五人格only          565 -> 566
$ws1.Range("F18").Value = 37   # 合肥·国乙only宇宙心动（含夜场） 36 -> 37

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = 168   # 合肥·Yolo Fes永乐庆典Vol.3 DAY1&DAY3 128 -> 168

# --- Sheet: 全部类型 (All types, combined) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = 168   # 合肥·Yolo Fes永乐庆典Vol.3 DAY1&DAY3 128 -> 168
$ws4.Range("F4").Value = 120   # 合肥·灵能百分百ONLY2.0      119 -> 120
$ws4.Range("F10").Value = 124  # 合肥·首届进击的巨人ONLY漫展  123 -> 124
$ws4.Range("F13").Value = 6899 # 合肥·第七届环形宇宙动漫游戏嘉年华 6891 -> 6899
$ws4.Range("F16").Value = 384  # 合肥·比翼连枝国乙&代号鸢only  383 -> 384
$ws4.Range("F17").Value = 3257 # 合肥·第八届环形宇宙动漫游戏嘉年华Plus 3244 -> 3257
$ws4.Range("F19").Value = 397  # ...水千丞签售预约票           395 -> 397
$ws4.Range("F21").Value = 566  # 合肥·SSS第五人格only          565 -> 566
$ws4.Range("F22").Value = 37   # 合肥·国乙only宇宙心动（含夜场） 36 -> 37
